# bioSample_hbrown_09.19.19.xlsx -- "continuing major accuracy cleaning"
#
# The floodmedia column (H2:H37) was recorded as numeric 0 placeholders;
# replace them with the explicit text "None" across all 36 data rows.
# Also refresh the view so the newly-edited column is selected/visible,
# and tighten the data-row height to match the rest of the cleanup pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data fix: H2:H37, numeric 0 -> text "None" ------------------
$ws.Range("H2:H37").Value = "None"

# --- Cosmetic: snug up the data row height (was 16) --------------------
$ws.Rows("2:37").RowHeight = 15

# --- Leave the cursor/selection on the column that was just edited -----
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H2:H37").Select()
